$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change %) scraped refresh.
# D-column price cells are forced back to Text (NumberFormat @ then
# ClearFormats) so numeric-looking strings like "63.43" are not silently
# coerced into real numbers by Excel's Range.Value auto-detection -- the
# source data stores these as plain text (t="inlineStr").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.711.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.592.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.30"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.05%  "
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.26"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.591.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.531"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.681.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.43"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "220.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.69"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("E24").Value = "  -3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.374.75"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("E34").Value = "  -4.98%  "
$ws.Range("E35").Value = "  -4.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.975"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.537"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.829"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("E46").Value = "  -5.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.728.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("E51").Value = "  -1.48%  "
